$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 432.05
$ws.Range("I19").Value = 357.44446
$ws.Range("J19").Value = 493.0909
$ws.Range("K19").Value = 357.44446
$ws.Range("L19").Value = 493.0909
$ws.Range("M19").Value = -182.44446
$ws.Range("N19").Value = -843.0908999999999
$ws.Range("H41").Value = 833.7826
$ws.Range("I41").Value = 673.3684
$ws.Range("K41").Value = 673.3684
$ws.Range("M41").Value = -233.3684
$ws.Range("H43").Value = 3944.6667
$ws.Range("J43").Value = 4000.3333
$ws.Range("L43").Value = 4000.3333
$ws.Range("N43").Value = -4138.3333
$ws.Range("H55").Value = 512.2105
$ws.Range("I55").Value = 532.375
$ws.Range("K55").Value = 532.375
$ws.Range("M55").Value = -318.375
$ws.Range("H87").Value = 28130.39
$ws.Range("J87").Value = 28130.39
$ws.Range("L87").Value = 28130.39
$ws.Range("N87").Value = -30626.39
$ws.Range("H90").Value = 28130.39
$ws.Range("J90").Value = 28130.39
$ws.Range("L90").Value = 84391.17
$ws.Range("N90").Value = -96871.17
$ws.Range("H92").Value = 1659
$ws.Range("I92").Value = 1264.3636
$ws.Range("K92").Value = 1264.3636
$ws.Range("M92").Value = -16.36359999999991
$ws.Range("H98").Value = 2710.353
$ws.Range("I98").Value = 3276.5386
$ws.Range("K98").Value = 3276.5386
$ws.Range("M98").Value = -1778.5386
$ws.Range("H116").Value = 17416.666
$ws.Range("J116").Value = 9700
$ws.Range("L116").Value = 9700
$ws.Range("N116").Value = -16584
$ws.Range("H122").Value = 2710.353
$ws.Range("I122").Value = 3276.5386
$ws.Range("K122").Value = 9829.6158
$ws.Range("M122").Value = -7379.6158
$ws.Range("H132").Value = 19136.69
$ws.Range("I132").Value = 1799.0385
$ws.Range("J132").Value = 169396.33
$ws.Range("K132").Value = 5397.1155
$ws.Range("L132").Value = 508188.99
$ws.Range("M132").Value = -2867.1155
$ws.Range("N132").Value = -513248.99
$ws.Range("H138").Value = 3498.7273
$ws.Range("I138").Value = 3198.6
$ws.Range("J138").Value = 6500
$ws.Range("K138").Value = 9595.799999999999
$ws.Range("L138").Value = 19500
$ws.Range("M138").Value = -4455.799999999999
$ws.Range("N138").Value = -29780

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2920.875
$ws.Range("I61").Value = 3066.7144
$ws.Range("K61").Value = 3066.7144
$ws.Range("M61").Value = -2854.7144
$ws.Range("H74").Value = 2022.3182
$ws.Range("I74").Value = 1700.4
$ws.Range("K74").Value = 1700.4
$ws.Range("M74").Value = -826.4000000000001
$ws.Range("H77").Value = 2022.3182
$ws.Range("I77").Value = 1700.4
$ws.Range("K77").Value = 8502
$ws.Range("M77").Value = -4134
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H97").Value = 2955.158
$ws.Range("I97").Value = 2236.6765
$ws.Range("J97").Value = 9062.25
$ws.Range("K97").Value = 2236.6765
$ws.Range("L97").Value = 9062.25
$ws.Range("M97").Value = -1740.6765
$ws.Range("N97").Value = -10054.25
$ws.Range("H132").Value = 1714.0769
$ws.Range("I132").Value = 1184.6
$ws.Range("J132").Value = 3479
$ws.Range("K132").Value = 3553.8
$ws.Range("L132").Value = 10437
$ws.Range("M132").Value = -1023.8
$ws.Range("N132").Value = -15497
$ws.Range("H135").Value = 51173.5
$ws.Range("J135").Value = 51173.5
$ws.Range("L135").Value = 51173.5
$ws.Range("N135").Value = -61313.5
$ws.Range("H136").Value = 2920.875
$ws.Range("I136").Value = 3066.7144
$ws.Range("K136").Value = 9200.143199999999
$ws.Range("M136").Value = -6650.143199999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 101853.7
$ws.Range("I20").Value = 2049.8572
$ws.Range("J20").Value = 334729.34
$ws.Range("K20").Value = 2049.8572
$ws.Range("L20").Value = 334729.34
$ws.Range("M20").Value = -1802.8572
$ws.Range("N20").Value = -335223.34
$ws.Range("H22").Value = 616.6667
$ws.Range("I22").Value = 616.6667
$ws.Range("K22").Value = 616.6667
$ws.Range("M22").Value = -443.6667
$ws.Range("H86").Value = 8046
$ws.Range("I86").Value = 7221
$ws.Range("K86").Value = 7221
$ws.Range("M86").Value = -6098
$ws.Range("H89").Value = 8046
$ws.Range("I89").Value = 7221
$ws.Range("K89").Value = 36105
$ws.Range("M89").Value = -30489
$ws.Range("H105").Value = 3445.6365
$ws.Range("I105").Value = 2819.625
$ws.Range("K105").Value = 2819.625
$ws.Range("M105").Value = -1072.625
$ws.Range("H107").Value = 5415.0835
$ws.Range("I107").Value = 4998.1
$ws.Range("K107").Value = 4998.1
$ws.Range("M107").Value = -3078.1
$ws.Range("H134").Value = 7944.8965
$ws.Range("I134").Value = 1806.1818
$ws.Range("J134").Value = 27238
$ws.Range("K134").Value = 5418.5454
$ws.Range("L134").Value = 81714
$ws.Range("M134").Value = -2883.5454
$ws.Range("N134").Value = -86784
$ws.Range("H135").Value = 57333.332
$ws.Range("J135").Value = 57333.332
$ws.Range("L135").Value = 57333.332
$ws.Range("N135").Value = -67473.33199999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3835
$ws.Range("I134").Value = 3483
$ws.Range("K134").Value = 10449
$ws.Range("M134").Value = -7914

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 790.1905
$ws.Range("J5").Value = 432.88235
$ws.Range("L5").Value = 1298.64705
$ws.Range("N5").Value = -1522.64705
$ws.Range("H18").Value = 731.75
$ws.Range("I18").Value = 731.75
$ws.Range("K18").Value = 2195.25
$ws.Range("M18").Value = -2026.25
$ws.Range("H68").Value = 1366.7567
$ws.Range("J68").Value = 1366.7567
$ws.Range("L68").Value = 4100.2701
$ws.Range("N68").Value = -5722.2701
$ws.Range("H71").Value = 1366.7567
$ws.Range("J71").Value = 1366.7567
$ws.Range("L71").Value = 12300.8103
$ws.Range("N71").Value = -20412.8103
$ws.Range("H134").Value = 5777.4443
$ws.Range("I134").Value = 2999.25
$ws.Range("K134").Value = 8997.75
$ws.Range("M134").Value = -3927.75
$ws.Range("H135").Value = 790.1905
$ws.Range("J135").Value = 432.88235
$ws.Range("L135").Value = 3895.94115
$ws.Range("N135").Value = -8965.941149999999
$ws.Range("H138").Value = 6122.154
$ws.Range("I138").Value = 3126.1667
$ws.Range("K138").Value = 9378.500100000001
$ws.Range("M138").Value = -4238.500100000001
$ws.Range("H139").Value = 62502400
$ws.Range("I139").Value = 76925130
$ws.Range("K139").Value = 230775390
$ws.Range("M139").Value = -230770250
$ws.Range("H141").Value = 500000420
$ws.Range("I141").Value = 500000420
$ws.Range("K141").Value = 1500001260
$ws.Range("M141").Value = -1499996080

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 185002000
$ws.Range("J80").Value = 185002000
$ws.Range("L80").Value = 185002000
$ws.Range("N80").Value = -185003996
$ws.Range("H83").Value = 185002000
$ws.Range("J83").Value = 185002000
$ws.Range("L83").Value = 925010000
$ws.Range("N83").Value = -925019984
$ws.Range("H93").Value = 70251
$ws.Range("J93").Value = 70251
$ws.Range("L93").Value = 70251
$ws.Range("N93").Value = -73995

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2496.9285
$ws.Range("I22").Value = 2224.5454
$ws.Range("K22").Value = 2224.5454
$ws.Range("M22").Value = -1929.5454
$ws.Range("H27").Value = 2496.9285
$ws.Range("I27").Value = 2224.5454
$ws.Range("K27").Value = 2224.5454
$ws.Range("M27").Value = -2117.5454
$ws.Range("H46").Value = 2634.7
$ws.Range("I46").Value = 1369.4
$ws.Range("K46").Value = 1369.4
$ws.Range("M46").Value = -1181.4
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751
$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992
$ws.Range("H136").Value = 6466.2593
$ws.Range("I136").Value = 5733.5835
$ws.Range("J136").Value = 7052.4
$ws.Range("K136").Value = 17200.7505
$ws.Range("L136").Value = 21157.2
$ws.Range("M136").Value = -14650.7505
$ws.Range("N136").Value = -26257.2
